$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.650.23'
$ws.Range('E2').Value = '  -2.28%  '
$ws.Range('D3').Value = '1.591.66'
$ws.Range('E3').Value = '  -2.50%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.00'
$ws.Range('E5').Value = '  -2.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.509'
$ws.Range('E6').Value = '  -2.35%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -1.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.66'
$ws.Range('E10').Value = '  -3.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0838'
$ws.Range('E11').Value = '  -1.45%  '
$ws.Range('D12').Value = '1.813.78'
$ws.Range('E12').Value = '  -2.50%  '
$ws.Range('D13').Value = '1.607.47'
$ws.Range('E13').Value = '  -1.24%  '
$ws.Range('E14').Value = '  -2.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.528'
$ws.Range('E15').Value = '  -3.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.74'
$ws.Range('E16').Value = '  -0.43%  '
$ws.Range('D17').Value = '26.666.77'
$ws.Range('E17').Value = '  -1.93%  '
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '207.85'
$ws.Range('E19').Value = '  -4.26%  '
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('E21').Value = '  -2.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.25'
$ws.Range('E22').Value = '  -2.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.37'
$ws.Range('E23').Value = '  -3.45%  '
$ws.Range('E24').Value = '  -2.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.20'
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.34'
$ws.Range('E27').Value = '  +0.76%  '
$ws.Range('E28').Value = '  -3.35%  '
$ws.Range('E29').Value = '  -2.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0504'
$ws.Range('E30').Value = '  -0.68%  '
$ws.Range('E31').Value = '  -2.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.26'
$ws.Range('E32').Value = '  -3.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.663'
$ws.Range('E33').Value = '  +22.57%  '
$ws.Range('D34').Value = '1.329.09'
$ws.Range('E34').Value = '  +0.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.91'
$ws.Range('E35').Value = '  -3.04%  '
$ws.Range('E36').Value = '  -3.19%  '
$ws.Range('E37').Value = '  -2.14%  '
$ws.Range('E38').Value = '  -1.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.828'
$ws.Range('E39').Value = '  -2.40%  '
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('E41').Value = '  +3.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.786'
$ws.Range('E42').Value = '  -1.64%  '
$ws.Range('E43').Value = '  -3.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.59'
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').Value = '1.726.59'
$ws.Range('E45').Value = '  -2.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.98'
$ws.Range('E46').Value = '  -0.91%  '
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.835'
$ws.Range('E48').Value = '  +2.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0510'
$ws.Range('E49').Value = '  -1.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0973'
$ws.Range('E50').Value = '  +1.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.52'
$ws.Range('E51').Value = '  -0.62%  '
